$d = $word.ActiveDocument

# Replace a whole paragraph's visible text with $newText while leaving any
# other runs in that paragraph (e.g. a leading empty "<w:r/>") untouched, and
# optionally applying run formatting (bold/italic) via $rPrXml (e.g. "<w:b/>").
# We scope the edit to the paragraph's own text run via Range.InsertXML so the
# zero-length sibling run survives (plain Range.Text / Find.Execute on a multi
# -run paragraph collapses/merges those sibling runs).
function Set-ParagraphText {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText,
        [string]$RPrXml = $null
    )

    $p = $d.Paragraphs($Index)
    $full = $p.Range
    # Exclude the trailing paragraph mark from the replaced range.
    $target = $d.Range($full.Start, $full.End - 1)

    if ($target.Text -ne $OldText) {
        throw "Paragraph $Index text mismatch: expected [$OldText] but found [$($target.Text)]"
    }

    $rPr = ""
    if ($RPrXml) { $rPr = "<w:rPr>$RPrXml</w:rPr>" }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $rPr + '<w:t>' + $NewText + '</w:t></w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

# Locate paragraphs by their current text so the script is resilient to the
# exact paragraph numbering.
function Find-ParagraphIndex([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $full = $p.Range
        $t = $d.Range($full.Start, $full.End - 1).Text
        if ($t -eq $text) { return $i }
    }
    throw "Paragraph not found: [$text]"
}

$oldTitle = "Play Mega Greatest Catch Free - Exciting Slot with Fishing and Dynamite Features"
$newTitle = "Play Mega Greatest Catch Free Slot Game"

# Page heading (Heading1 style).
Set-ParagraphText (Find-ParagraphIndex $oldTitle) $oldTitle $newTitle

# "What we like" bullet list.
Set-ParagraphText (Find-ParagraphIndex "High-definition graphics with relaxing sound.") `
    "High-definition graphics with relaxing sound." "High-definition graphics"

Set-ParagraphText (Find-ParagraphIndex "Fishing and Dynamite features offer exciting gameplay.") `
    "Fishing and Dynamite features offer exciting gameplay." "Relaxing sound"

Set-ParagraphText (Find-ParagraphIndex "Autoplay and Turbo functions for faster play.") `
    "Autoplay and Turbo functions for faster play." "Exciting bonus features"

Set-ParagraphText (Find-ParagraphIndex "Medium volatility for a wide range of players.") `
    "Medium volatility for a wide range of players." "Wide compatibility"

# "What we don't like" bullet list.
Set-ParagraphText (Find-ParagraphIndex "Only 12 fixed paylines.") `
    "Only 12 fixed paylines." "Limited number of paylines"

$oldBetText = "Minimum bet starts at " + [char]0x20AC + "0.10, which may be too high for some players."
Set-ParagraphText (Find-ParagraphIndex $oldBetText) $oldBetText "Limited betting range"

# Bold "play now" call-to-action line (repeats the page heading, bold).
Set-ParagraphText (Find-ParagraphIndex $oldTitle) $oldTitle $newTitle "<w:b/>"

# Italic meta-description line.
$oldMeta = "Read our review of Mega Greatest Catch, an HTML 5 online slot with medium volatility, high-definition graphics, and fishing and dynamite features. Play now for free."
$newMeta = "Read our review of Mega Greatest Catch and play this exciting slot game for free."
Set-ParagraphText (Find-ParagraphIndex $oldMeta) $oldMeta $newMeta "<w:i/>"
